# Update the "取得日時" (retrieved-at) timestamp in column A for every
# data row on the "ランサーズ" sheet to reflect the new run time:
#   2026-01-15 01:25:52  ->  2026-01-15 01:59:18

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-15 01:59:18"

# Data rows are 2 through 18 (row 1 is the header).
for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
